$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recomputed values in row 2
$ws.Range("AB2").Value = 64.37435589208292
$ws.Range("AC2").Value = 3.949086466552615
$ws.Range("AD2").Value = 13.45632043806779
$ws.Range("AG2").Value = 81.428519079832
$ws.Range("AH2").Value = 0.572542084751227

# Delete the row for sample UK026 (row 3), shifting cells up
$ws.Rows("3").Delete()
